## TC47_Canine_Filter_Breed-YorkshireTerr.xlsx update
## - Reformats the "StatQuery" (column C) Cypher query on all three tab rows
##   (CasesTab / SamplesTab / FilesTab) from a single compact line to a
##   nicely indented multi-line version (splitting the combined OPTIONAL
##   MATCH clause into two statements).
## - Moves the active selection from C2 to D2 and scrolls the view so row 2
##   is the top visible row.
## - Row 3's height grows (217.5 -> 246.5) to fit the taller wrapped text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newStatQuery = @"
MATCH (s:study)
  WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies
  MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies
  MATCH (d:diagnosis)
  WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies
  MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
    WHERE demo.breed IN ['Yorkshire Terrier']
  OPTIONAL MATCH (f:file)-[*]->(c)
  OPTIONAL MATCH (samp:sample)-[*]->(c)
  WITH DISTINCT c AS c, p, s, demo, diag, f, samp
  RETURN count(DISTINCT(f)) as number_of_files ,
             count(DISTINCT(samp)) as number_of_sample ,
             count(DISTINCT(c.case_id)) as number_of_cases ,
             count(DISTINCT(s.clinical_study_designation)) as number_of_study
"@

# Same reformatted StatQuery text lands in C2, C3 and C4 - the per-tab
# queries in column B are untouched.
$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

# Row 3 grows to fit the now-taller wrapped StatQuery text.
$ws.Rows.Item(3).RowHeight = 246.5

# Selection moves from C2 to D2, with the view scrolled so row 2 is on top.
$ws.Range("D2").Select()

Write-Output "edit applied"
